# Update "想去人数" (want-to-go count) values in column F across sheets
# "展览", "演出" and "全部类型", reflecting a freshly regenerated scrape
# of the source data (gh-pages output generated at 456a3b4).

$wb = $excel.ActiveWorkbook

# Sheet "展览"
$ws1 = $wb.Worksheets.Item("展览")
$ws1Updates = @{
    2  = 2040
    6  = 69
    8  = 370
    11 = 944
    12 = 301
    16 = 355
    17 = 324
    18 = 738
    20 = 695
    21 = 237
    23 = 954
    24 = 413
    27 = 340
    29 = 31
    30 = 444
}
foreach ($row in $ws1Updates.Keys) {
    $ws1.Range("F$row").Value = $ws1Updates[$row]
}

# Sheet "演出"
$ws2 = $wb.Worksheets.Item("演出")
$ws2Updates = @{
    5  = 29
    12 = 31
}
foreach ($row in $ws2Updates.Keys) {
    $ws2.Range("F$row").Value = $ws2Updates[$row]
}

# Sheet "全部类型"
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4Updates = @{
    3  = 2040
    8  = 69
    10 = 370
    13 = 944
    14 = 301
    20 = 29
    21 = 355
    24 = 324
    25 = 738
    27 = 695
    28 = 237
    30 = 954
    31 = 413
    36 = 340
    40 = 31
    41 = 31
    42 = 444
}
foreach ($row in $ws4Updates.Keys) {
    $ws4.Range("F$row").Value = $ws4Updates[$row]
}
